$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected Price/Volume columns to Text format so numeric-looking
# strings (e.g. "27.552.27", "0.9988", "  -0.87%  ") are preserved exactly as text,
# matching the original inline-string cell contents, then strip the format
# override afterwards so no visible style change is introduced.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.552.27"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.754.16"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "324.64"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "0.4464"
$ws.Range("E7").Value = "  +5.10%  "
$ws.Range("D8").Value = "0.3583"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "0.07501"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "41.96"
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").Value = "1.095"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "20.82"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "6.029"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "7.119"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "1.751.38"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").Value = "93.29"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "0.00001062"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "0.06411"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "0.9986"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "16.84"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "5.817"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").Value = "27.596.22"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "2.106"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("D26").Value = "162.83"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "20.48"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").Value = "1.953.95"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").Value = "2.092"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "126.24"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "1.078"
$ws.Range("E31").Value = "  -7.62%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.09080"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "3.666"
$ws.Range("E33").Value = "  +4.29%  "
$ws.Range("D34").Value = "5.536"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("D36").Value = "0.02291"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "0.2100"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "0.06032"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "0.6368"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "4.962"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "1.203"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").Value = "1.375"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").Value = "7.800"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "13.31"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "0.5915"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "3.712"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "122.43"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").Value = "1.954"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "0.06853"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "72.52"
$ws.Range("E51").Value = "  -2.44%  "

$ws.Range("D2:E51").ClearFormats()
